$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.917.19'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.640.58'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5034'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.70%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.004'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2569'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06398'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.48%  '
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07731'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.269'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.641.63'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.862.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5463'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅7923'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.45'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.914.73'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '203.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.382'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.927'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.990'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.931'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.44'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.99%  '
$ws.Range('E27').Value = '  -3.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.69'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.737'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.245'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04939'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.280'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.192'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.546'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.380'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.93%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.637'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.68%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.8957'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.162.96'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5619'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01565'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.004'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.729'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8100'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.773.58'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈117'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4531'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05059'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.004'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.33%  '
